# The "Recorded By" column (G) lists the names/emails of whoever logged the
# attendance for a session. Many rows show "dnasr281@gmail.com, System" -
# this commit reorders that list so "System" is listed first, i.e.
# "System, dnasr281@gmail.com".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
